$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-26 (columns A-I) with the new metric values / reordered model names.

$ws.Range("A2").Value = "model_7_1_0"
$ws.Range("B2").Value = 0.1758715878431283
$ws.Range("C2").Value = 0.2165196650392641
$ws.Range("D2").Value = -2.594715022818576
$ws.Range("E2").Value = -0.09679580481713934
$ws.Range("F2").Value = 0.9120672941207886
$ws.Range("G2").Value = 1.123998284339905
$ws.Range("H2").Value = 0.7450722455978394
$ws.Range("I2").Value = 0.945681631565094
$ws.Range("A3").Value = "model_7_1_1"
$ws.Range("B3").Value = 0.2476815344206118
$ws.Range("C3").Value = 0.1163988428477337
$ws.Range("D3").Value = -3.459487575920075
$ws.Range("E3").Value = -0.2828150818619202
$ws.Range("F3").Value = 0.8325948715209961
$ws.Range("G3").Value = 1.26763391494751
$ws.Range("H3").Value = 0.9243125915527344
$ws.Range("I3").Value = 1.106071591377258
$ws.Range("A4").Value = "model_7_1_2"
$ws.Range("B4").Value = 0.2994798195068017
$ws.Range("C4").Value = 0.1622551767936917
$ws.Range("D4").Value = -2.600093777275944
$ws.Range("E4").Value = -0.145201434254318
$ws.Range("F4").Value = 0.7752693891525269
$ws.Range("G4").Value = 1.201847434043884
$ws.Range("H4").Value = 0.7461870908737183
$ws.Range("I4").Value = 0.9874181151390076
$ws.Range("A5").Value = "model_7_1_3"
$ws.Range("B5").Value = 0.3084822613418926
$ws.Range("C5").Value = 0.1962681791300993
$ws.Range("D5").Value = -3.105354642684081
$ws.Range("E5").Value = -0.1723990891733698
$ws.Range("F5").Value = 0.7653063535690308
$ws.Range("G5").Value = 1.153051614761353
$ws.Range("H5").Value = 0.8509118556976318
$ws.Range("I5").Value = 1.010868549346924
$ws.Range("A6").Value = "model_7_1_4"
$ws.Range("B6").Value = 0.3196790766247024
$ws.Range("C6").Value = 0.2339345811183091
$ws.Range("D6").Value = -3.574437695792332
$ws.Range("E6").Value = -0.1922818181893866
$ws.Range("F6").Value = 0.7529147863388062
$ws.Range("G6").Value = 1.099014401435852
$ws.Range("H6").Value = 0.9481382369995117
$ws.Range("I6").Value = 1.028011798858643
$ws.Range("A7").Value = "model_7_1_5"
$ws.Range("B7").Value = 0.5517609040528361
$ws.Range("C7").Value = 0.4374461518999192
$ws.Range("D7").Value = -2.181865917092365
$ws.Range("E7").Value = 0.1445147251889901
$ws.Range("F7").Value = 0.4960685670375824
$ws.Range("G7").Value = 0.8070523142814636
$ws.Range("H7").Value = 0.6595014929771423
$ws.Range("I7").Value = 0.7376183271408081
$ws.Range("A8").Value = "model_7_1_6"
$ws.Range("B8").Value = 0.699451206184793
$ws.Range("C8").Value = 0.666415645528635
$ws.Range("D8").Value = -1.874724129541447
$ws.Range("E8").Value = 0.3809558300358545
$ws.Range("F8").Value = 0.3326189517974854
$ws.Range("G8").Value = 0.4785675406455994
$ws.Range("H8").Value = 0.5958405733108521
$ws.Range("I8").Value = 0.5337535738945007
$ws.Range("A9").Value = "model_7_1_7"
$ws.Range("B9").Value = 0.7202407469124898
$ws.Range("C9").Value = 0.6837177400934029
$ws.Range("D9").Value = -1.465726588896493
$ws.Range("E9").Value = 0.4424642220398964
$ws.Range("F9").Value = 0.309611052274704
$ws.Range("G9").Value = 0.4537455141544342
$ws.Range("H9").Value = 0.5110681056976318
$ws.Range("I9").Value = 0.4807197451591492
$ws.Range("A10").Value = "model_7_1_9"
$ws.Range("B10").Value = 0.7233798777770134
$ws.Range("C10").Value = 0.6783971650634322
$ws.Range("D10").Value = -1.652391395273584
$ws.Range("E10").Value = 0.4166569037917317
$ws.Range("F10").Value = 0.3061369359493256
$ws.Range("G10").Value = 0.4613785743713379
$ws.Range("H10").Value = 0.5497579574584961
$ws.Range("I10").Value = 0.5029712915420532
$ws.Range("A11").Value = "model_7_1_8"
$ws.Range("B11").Value = 0.7283382578375334
$ws.Range("C11").Value = 0.6894314283802769
$ws.Range("D11").Value = -1.420979368097305
$ws.Range("E11").Value = 0.4525589037214911
$ws.Range("F11").Value = 0.300649493932724
$ws.Range("G11").Value = 0.4455485939979553
$ws.Range("H11").Value = 0.5017934441566467
$ws.Range("I11").Value = 0.4720158874988556
$ws.Range("A12").Value = "model_7_1_10"
$ws.Range("B12").Value = 0.7412982294612298
$ws.Range("C12").Value = 0.6874241898322258
$ws.Range("D12").Value = -1.526178561057334
$ws.Range("E12").Value = 0.4388899010280761
$ws.Range("F12").Value = 0.2863066494464874
$ws.Range("G12").Value = 0.448428213596344
$ws.Range("H12").Value = 0.5235980153083801
$ws.Range("I12").Value = 0.4838015735149384
$ws.Range("A13").Value = "model_7_1_11"
$ws.Range("B13").Value = 0.7490160177648111
$ws.Range("C13").Value = 0.690416243979882
$ws.Range("D13").Value = -1.44670795015799
$ws.Range("E13").Value = 0.4505123795373531
$ws.Range("F13").Value = 0.277765303850174
$ws.Range("G13").Value = 0.4441357553005219
$ws.Range("H13").Value = 0.5071262121200562
$ws.Range("I13").Value = 0.4737803936004639
$ws.Range("A14").Value = "model_7_1_12"
$ws.Range("B14").Value = 0.7558286995599944
$ws.Range("C14").Value = 0.6909260225215247
$ws.Range("D14").Value = -1.370876859501327
$ws.Range("E14").Value = 0.4595436448186199
$ws.Range("F14").Value = 0.2702256739139557
$ws.Range("G14").Value = 0.443404346704483
$ws.Range("H14").Value = 0.4914087653160095
$ws.Range("I14").Value = 0.4659934341907501
$ws.Range("A15").Value = "model_7_1_13"
$ws.Range("B15").Value = 0.7600408125903589
$ws.Range("C15").Value = 0.6909714015954092
$ws.Range("D15").Value = -1.381159726057862
$ws.Range("E15").Value = 0.458416805738438
$ws.Range("F15").Value = 0.2655641138553619
$ws.Range("G15").Value = 0.4433393180370331
$ws.Range("H15").Value = 0.4935401082038879
$ws.Range("I15").Value = 0.4669650197029114
$ws.Range("A16").Value = "model_7_1_14"
$ws.Range("B16").Value = 0.7642816107488581
$ws.Range("C16").Value = 0.6908593613426857
$ws.Range("D16").Value = -1.391423024256448
$ws.Range("E16").Value = 0.4571563145747096
$ws.Range("F16").Value = 0.2608707845211029
$ws.Range("G16").Value = 0.4435000419616699
$ws.Range("H16").Value = 0.4956673681735992
$ws.Range("I16").Value = 0.4680518507957458
$ws.Range("A17").Value = "model_7_1_15"
$ws.Range("B17").Value = 0.7684399363964868
$ws.Range("C17").Value = 0.6896039299916821
$ws.Range("D17").Value = -1.38121254152776
$ws.Range("E17").Value = 0.4572057184080579
$ws.Range("F17").Value = 0.2562687695026398
$ws.Range("G17").Value = 0.4453011155128479
$ws.Range("H17").Value = 0.4935510456562042
$ws.Range("I17").Value = 0.4680092334747314
$ws.Range("A18").Value = "model_7_1_16"
$ws.Range("B18").Value = 0.7722367770099429
$ws.Range("C18").Value = 0.6881817621065127
$ws.Range("D18").Value = -1.378049442278153
$ws.Range("E18").Value = 0.4563127488399654
$ws.Range("F18").Value = 0.2520667910575867
$ws.Range("G18").Value = 0.4473413825035095
$ws.Range("H18").Value = 0.4928954243659973
$ws.Range("I18").Value = 0.4687792062759399
$ws.Range("A19").Value = "model_7_1_17"
$ws.Range("B19").Value = 0.774730422510727
$ws.Range("C19").Value = 0.6877957477703214
$ws.Range("D19").Value = -1.425406646738312
$ws.Range("E19").Value = 0.4506135310620867
$ws.Range("F19").Value = 0.2493070214986801
$ws.Range("G19").Value = 0.4478951990604401
$ws.Range("H19").Value = 0.5027111172676086
$ws.Range("I19").Value = 0.4736931622028351
$ws.Range("A20").Value = "model_7_1_18"
$ws.Range("B20").Value = 0.7781328911552264
$ws.Range("C20").Value = 0.6831967981947833
$ws.Range("D20").Value = -1.423943180695787
$ws.Range("E20").Value = 0.4467310151097583
$ws.Range("F20").Value = 0.2455415278673172
$ws.Range("G20").Value = 0.4544928967952728
$ws.Range("H20").Value = 0.5024077892303467
$ws.Range("I20").Value = 0.4770407378673553
$ws.Range("A21").Value = "model_7_1_24"
$ws.Range("B21").Value = 0.7882515426542933
$ws.Range("C21").Value = 0.6587002796008148
$ws.Range("D21").Value = -2.165987982074765
$ws.Range("E21").Value = 0.3412077838699829
$ws.Range("F21").Value = 0.2343431562185287
$ws.Range("G21").Value = 0.4896361529827118
$ws.Range("H21").Value = 0.6562104821205139
$ws.Range("I21").Value = 0.5680251717567444
$ws.Range("A22").Value = "model_7_1_21"
$ws.Range("B22").Value = 0.7888522598718374
$ws.Range("C22").Value = 0.6527966613350394
$ws.Range("D22").Value = -1.825157191147835
$ws.Range("E22").Value = 0.3745646571549697
$ws.Range("F22").Value = 0.233678326010704
$ws.Range("G22").Value = 0.4981056451797485
$ws.Range("H22").Value = 0.5855668783187866
$ws.Range("I22").Value = 0.5392642021179199
$ws.Range("A23").Value = "model_7_1_22"
$ws.Range("B23").Value = 0.7891045914586325
$ws.Range("C23").Value = 0.6576895427300058
$ws.Range("D23").Value = -1.94530636732738
$ws.Range("E23").Value = 0.3652840927267299
$ws.Range("F23").Value = 0.2333990484476089
$ws.Range("G23").Value = 0.4910861551761627
$ws.Range("H23").Value = 0.6104701161384583
$ws.Range("I23").Value = 0.5472661256790161
$ws.Range("A24").Value = "model_7_1_23"
$ws.Range("B24").Value = 0.7892802020620986
$ws.Range("C24").Value = 0.6579264712615286
$ws.Range("D24").Value = -2.043191474977315
$ws.Range("E24").Value = 0.3544198338952598
$ws.Range("F24").Value = 0.2332047075033188
$ws.Range("G24").Value = 0.4907462894916534
$ws.Range("H24").Value = 0.6307585835456848
$ws.Range("I24").Value = 0.5566335320472717
$ws.Range("A25").Value = "model_7_1_20"
$ws.Range("B25").Value = 0.7894266102101462
$ws.Range("C25").Value = 0.6536621707324612
$ws.Range("D25").Value = -1.703826178998256
$ws.Range("E25").Value = 0.3890513552449202
$ws.Range("F25").Value = 0.2330427020788193
$ws.Range("G25").Value = 0.4968639612197876
$ws.Range("H25").Value = 0.5604187250137329
$ws.Range("I25").Value = 0.5267734527587891
$ws.Range("A26").Value = "model_7_1_19"
$ws.Range("B26").Value = 0.7900246628342548
$ws.Range("C26").Value = 0.6502201847162494
$ws.Range("D26").Value = -1.555407595326259
$ws.Range("E26").Value = 0.4028095126525227
$ws.Range("F26").Value = 0.2323808372020721
$ws.Range("G26").Value = 0.5018019080162048
$ws.Range("H26").Value = 0.5296562314033508
$ws.Range("I26").Value = 0.514910876750946
